$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 242.48983023918456
$ws.Range("D2").Value = 201.43216889294126
$ws.Range("E2").Value = 293.7617597245353

$ws.Range("C3").Value = 294.88681562904264
$ws.Range("D3").Value = 246.86629098205836
$ws.Range("E3").Value = 354.225585119166

$ws.Range("C4").Value = 281.46618710131924
$ws.Range("D4").Value = 240.22707503781177
$ws.Range("E4").Value = 331.26051594912161
$ws.Range("F4").Value = 90
$ws.Range("G4").Value = 281.62545755195231
$ws.Range("H4").Value = 249.20151924932676
$ws.Range("I4").Value = 25.492657587262556
$ws.Range("J4").Value = 298.21002153687783
$ws.Range("K4").Value = 172.972731633774
$ws.Range("L4").Value = 501.32010073781998
$ws.Range("M4").Value = 235.2661414204735
$ws.Range("N4").Value = 5.8232862115540032

$ws.Range("C5").Value = 311.35914964864651
$ws.Range("D5").Value = 262.55991930463489
$ws.Range("E5").Value = 371.11514253201273
$ws.Range("F5").Value = 90
$ws.Range("G5").Value = 313.00084438763133
$ws.Range("H5").Value = 277.35192189145215
$ws.Range("I5").Value = 28.36957217583722
$ws.Range("J5").Value = 339.73112679327892
$ws.Range("K5").Value = 196.59419738297501
$ws.Range("L5").Value = 556.15024332354881
$ws.Range("M5").Value = 261.41722776113909
$ws.Range("N5").Value = 5.9288556405821646

$ws.Range("C6").Value = 370.2753152062495
$ws.Range("D6").Value = 308.01263296182378
$ws.Range("E6").Value = 447.71258598500236
$ws.Range("F6").Value = 90
$ws.Range("G6").Value = 372.68015284974774
$ws.Range("H6").Value = 347.37323743001178
$ws.Range("I6").Value = 35.42496277413386
$ws.Range("J6").Value = 433.87208090914334
$ws.Range("K6").Value = 256.33666210491504
$ws.Range("L6").Value = 680.55502750315281
$ws.Range("M6").Value = 308.54536639564083
$ws.Range("N6").Value = 6.1116957894524182

$ws.Range("C7").Value = 451.43783176235161
$ws.Range("D7").Value = 380.15020668433942
$ws.Range("E7").Value = 538.72167771060003
$ws.Range("F7").Value = 90
$ws.Range("G7").Value = 489.84795559987674
$ws.Range("H7").Value = 358.02616528808846
$ws.Range("I7").Value = 36.957831859392854
$ws.Range("J7").Value = 538.89668446521807
$ws.Range("K7").Value = 287.79837666053191
$ws.Range("L7").Value = 823.75523332093894
$ws.Range("M7").Value = 421.65762743009606
$ws.Range("N7").Value = 6.3457063743255082
